$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Rows("2:2").ClearContents()
$ws.Range("A2").Value = "EMEA"
$ws.Range("B2").Value = "BCPT LN"
$ws.Range("C2").Value = "Balanced Comm Ppty Trust"
$ws.Range("D2").Value = 45610
$ws.Range("E2").Value = "E"
$ws.Range("F2").Value = "MSCI DM Small Cap"
$ws.Range("G2").Value = "Delete"
$ws.Range("H2").Value = -0.0001081784937401948
$ws.Range("I2").Value = -17.40727434825358
$ws.Range("J2").Value = -14.03100374
$ws.Range("K2").Value = -3.819823859640908
$ws.Range("L2").Value = -17.40727434825358
$ws.Range("M2").Value = -14.03100374
$ws.Range("N2").Value = -3.819823859640908
$ws.Range("O2").Value = 45595
$ws.Range("P2").Value = "Acquisition"
$ws.Range("Q2").Value = 45604
$ws.Range("W2").Value = "TRUE"

# Row 3
$ws.Rows("3:3").ClearContents()
$ws.Range("A3").Value = "EMEA"
$ws.Range("B3").Value = "ANG SJ"
$ws.Range("C3").Value = "Anglogold Ashanti"
$ws.Range("D3").Value = 45617
$ws.Range("E3").Value = "E"
$ws.Range("F3").Value = "MSCI EM"
$ws.Range("G3").Value = "S Inc"
$ws.Range("H3").Value = 0.0002609298317320255
$ws.Range("I3").Value = 169.2567887469962
$ws.Range("J3").Value = 6.517444245722097
$ws.Range("K3").Value = 3.471595128764175
$ws.Range("L3").Value = 169.2567887469962
$ws.Range("M3").Value = 6.517444245722097
$ws.Range("N3").Value = 3.471595128764175
$ws.Range("O3").Value = 45602
$ws.Range("P3").Value = "Acquisition"
$ws.Range("Q3").Value = 45604
$ws.Range("R3").Value = "TRUE"

# Row 4
$ws.Rows("4:4").ClearContents()
$ws.Range("A4").Value = "EMEA"
$ws.Range("B4").Value = "CEY LN"
$ws.Range("C4").Value = "Centamin"
$ws.Range("D4").Value = 45617
$ws.Range("E4").Value = "E"
$ws.Range("F4").Value = "MSCI DM Small Cap"
$ws.Range("G4").Value = "Delete"
$ws.Range("H4").Value = -0.000278655968179756
$ws.Range("I4").Value = -45.0171034824237
$ws.Range("J4").Value = -23.2216539
$ws.Range("K4").Value = -1.690132056726749
$ws.Range("L4").Value = -45.0171034824237
$ws.Range("M4").Value = -23.2216539
$ws.Range("N4").Value = -1.690132056726749
$ws.Range("O4").Value = 45602
$ws.Range("P4").Value = "Acquisition"
$ws.Range("Q4").Value = 45604
$ws.Range("W4").Value = "TRUE"

# Row 5
$ws.Rows("5:5").ClearContents()
$ws.Range("A5").Value = "EMEA"
$ws.Range("B5").Value = "TBD"
$ws.Range("C5").Value = " SUNRISE A"
$ws.Range("D5").Value = 45621
$ws.Range("E5").Value = "E"
$ws.Range("F5").Value = "MSCI DM Small Cap"
$ws.Range("G5").Value = "Add"
$ws.Range("J5").Value = 1.1619372024
$ws.Range("M5").Value = 0.2977813214999999
$ws.Range("O5").Value = 45596
$ws.Range("P5").Value = "Spin-off from Liberty Global"
$ws.Range("Q5").Value = 45604
$ws.Range("V5").Value = "TRUE"
$ws.Range("X5").Value = "TRUE"

# Row 6
$ws.Rows("6:6").ClearContents()
$ws.Range("A6").Value = "EMEA"
$ws.Range("B6").Value = "TBD"
$ws.Range("C6").Value = " SUNRISE A"
$ws.Range("D6").Value = 45621
$ws.Range("E6").Value = "E"
$ws.Range("F6").Value = "MSCI DM Small Cap"
$ws.Range("G6").Value = "Spin-off Delete"
$ws.Range("J6").Value = -0.8641558809
$ws.Range("M6").Value = 0.2977813214999999
$ws.Range("O6").Value = 45596
$ws.Range("P6").Value = "Spin-off from Liberty Global"
$ws.Range("Q6").Value = 45604
$ws.Range("W6").Value = "TRUE"
$ws.Range("X6").Value = "TRUE"

# Row 7
$ws.Rows("7:7").ClearContents()
$ws.Range("A7").Value = "AP"
$ws.Range("B7").Value = "TBD"
$ws.Range("C7").Value = "Jusung Engineering (New)"
$ws.Range("D7").Value = 45631
$ws.Range("E7").Value = "C"
$ws.Range("F7").Value = "FTSE DM Small Cap"
$ws.Range("G7").Value = "Add"
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 0
$ws.Range("O7").Value = 45590
$ws.Range("P7").Value = "Demerger"
$ws.Range("Q7").Value = 45594
$ws.Range("V7").Value = "TRUE"

# Row 8
$ws.Rows("8:8").ClearContents()
$ws.Range("A8").Value = "AP"
$ws.Range("B8").Value = "TBD"
$ws.Range("C8").Value = "Jusung Engineering (New)"
$ws.Range("D8").Value = 45632
$ws.Range("E8").Value = "C"
$ws.Range("F8").Value = "MSCI EM Small Cap"
$ws.Range("G8").Value = "Add"
$ws.Range("M8").Value = -0.1196923767245363
$ws.Range("O8").Value = 45580
$ws.Range("P8").Value = "Spin-Off"
$ws.Range("Q8").Value = 45595
$ws.Range("V8").Value = "TRUE"
$ws.Range("X8").Value = "TRUE"

# Row 9
$ws.Rows("9:9").ClearContents()
$ws.Range("A9").Value = "AP"
$ws.Range("B9").Value = "TBD"
$ws.Range("C9").Value = "Jusung Engineering (New)"
$ws.Range("D9").Value = 45632
$ws.Range("E9").Value = "C"
$ws.Range("F9").Value = "KOSDAQ150"
$ws.Range("G9").Value = "Delete"
$ws.Range("J9").Value = -0.1196923767245363
$ws.Range("M9").Value = -0.1196923767245363
$ws.Range("O9").Value = 45588
$ws.Range("P9").Value = "Spin-Off by Jusung Engineering"
$ws.Range("Q9").Value = 45595
$ws.Range("W9").Value = "TRUE"
$ws.Range("X9").Value = "TRUE"

# Row 10
$ws.Rows("10:10").ClearContents()
$ws.Range("A10").Value = "EMEA"
$ws.Range("B10").Value = "TBD"
$ws.Range("C10").Value = "Canal+"
$ws.Range("D10").Value = 45642
$ws.Range("E10").Value = "C"
$ws.Range("F10").Value = "CAC 40"
$ws.Range("G10").Value = "Spin-off Delete"
$ws.Range("J10").Value = -12.976968375
$ws.Range("M10").Value = -12.976968375
$ws.Range("O10").Value = 45600
$ws.Range("P10").Value = "Spin-off from Vivendi"
$ws.Range("Q10").Value = 45604
$ws.Range("W10").Value = "TRUE"
$ws.Range("X10").Value = "TRUE"

# Row 11
$ws.Rows("11:11").ClearContents()
$ws.Range("A11").Value = "EMEA"
$ws.Range("B11").Value = "TBD"
$ws.Range("C11").Value = "Louis Hachette"
$ws.Range("D11").Value = 45642
$ws.Range("E11").Value = "C"
$ws.Range("F11").Value = "CAC 40"
$ws.Range("G11").Value = "Spin-off Delete"
$ws.Range("J11").Value = -12.976968375
$ws.Range("M11").Value = -12.976968375
$ws.Range("O11").Value = 45600
$ws.Range("P11").Value = "Spin-off from Vivendi"
$ws.Range("Q11").Value = 45604
$ws.Range("W11").Value = "TRUE"
$ws.Range("X11").Value = "TRUE"

# Row 12
$ws.Rows("12:12").ClearContents()
$ws.Range("A12").Value = "EMEA"
$ws.Range("B12").Value = "TBD"
$ws.Range("C12").Value = "Havas"
$ws.Range("D12").Value = 45642
$ws.Range("E12").Value = "C"
$ws.Range("F12").Value = "CAC 40"
$ws.Range("G12").Value = "Spin-off Delete"
$ws.Range("J12").Value = -12.976968375
$ws.Range("M12").Value = -12.976968375
$ws.Range("O12").Value = 45600
$ws.Range("P12").Value = "Spin-off from Vivendi"
$ws.Range("Q12").Value = 45604
$ws.Range("W12").Value = "TRUE"
$ws.Range("X12").Value = "TRUE"

# Delete rows 13 and 14 (shifted content already rewritten above)
$ws.Rows("13:14").Delete()
